$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

$ws.Range("G1").Value = "Quantity"
$ws.Range("G2").Value = 1

$ws.Range("G2").Select()
